# Fruta / hortaliza, semanal
#
# Weekly refresh of the "Feria Lagunitas de Puerto Montt - Pomelo" price
# series: the oldest weekly record (row 132, 2021-04-30) drops off, every
# remaining record (rows 133-168) shifts up one row, and a brand-new
# record (2022-03-16) is appended at the former last row (168).
#
# Columns D (Fecha), K (Variedad), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado) and
# S (Precio $/Kg) are rewritten per row; all other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 132; D = 44582; K = "Start Ruby"; L = "Primera"; M = 200; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 133; D = 44159; K = "Start Ruby"; L = "Primera"; M = 160; N = 9000; O = 10000; P = 9500; S = 679 },
    @{ Row = 134; D = 44386; K = "Red Blush"; L = "Primera"; M = 160; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 135; D = 44466; K = "Start Ruby"; L = "Primera"; M = 60; N = 12000; O = 12000; P = 12000; S = 857 },
    @{ Row = 136; D = 44322; K = "Start Ruby"; L = "Primera"; M = 120; N = 19000; O = 20000; P = 19500; S = 1393 },
    @{ Row = 137; D = 44452; K = "Start Ruby"; L = "Primera"; M = 40; N = 13000; O = 13000; P = 13000; S = 929 },
    @{ Row = 138; D = 44358; K = "Start Ruby"; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 139; D = 44460; K = "Start Ruby"; L = "Primera"; M = 80; N = 12000; O = 12000; P = 12000; S = 857 },
    @{ Row = 140; D = 44333; K = "Start Ruby"; L = "Primera"; M = 20; N = 19000; O = 20000; P = 19500; S = 1393 },
    @{ Row = 141; D = 44244; K = "Start Ruby"; L = "Primera"; M = 40; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 142; D = 44356; K = "Start Ruby"; L = "Primera"; M = 60; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 143; D = 44539; K = "Start Ruby"; L = "Primera"; M = 200; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 144; D = 44547; K = "Start Ruby"; L = "Primera"; M = 300; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 145; D = 44482; K = "Start Ruby"; L = "Primera"; M = 80; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 146; D = 44603; K = "Start Ruby"; L = "Primera"; M = 200; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 147; D = 44603; K = "Start Ruby"; L = "Segunda"; M = 100; N = 11000; O = 11000; P = 11000; S = 786 },
    @{ Row = 148; D = 44217; K = "Start Ruby"; L = "Primera"; M = 80; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 149; D = 44596; K = "Start Ruby"; L = "Primera"; M = 180; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 150; D = 44596; K = "Start Ruby"; L = "Segunda"; M = 60; N = 11000; O = 11000; P = 11000; S = 786 },
    @{ Row = 151; D = 44326; K = "Start Ruby"; L = "Primera"; M = 60; N = 19000; O = 20000; P = 19500; S = 1393 },
    @{ Row = 152; D = 44238; K = "Start Ruby"; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 153; D = 44242; K = "Start Ruby"; L = "Primera"; M = 60; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 154; D = 44166; K = "Start Ruby"; L = "Primera"; M = 200; N = 9000; O = 10000; P = 9500; S = 679 },
    @{ Row = 155; D = 44348; K = "Start Ruby"; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 156; D = 44223; K = "Start Ruby"; L = "Primera"; M = 60; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 157; D = 44579; K = "Start Ruby"; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 158; D = 44579; K = "Start Ruby"; L = "Segunda"; M = 100; N = 10000; O = 10000; P = 10000; S = 714 },
    @{ Row = 159; D = 44515; K = "Start Ruby"; L = "Primera"; M = 60; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 160; D = 44438; K = "Start Ruby"; L = "Primera"; M = 60; N = 11000; O = 11000; P = 11000; S = 786 },
    @{ Row = 161; D = 44249; K = "Start Ruby"; L = "Primera"; M = 60; N = 13000; O = 14000; P = 13500; S = 964 },
    @{ Row = 162; D = 44566; K = "Start Ruby"; L = "Primera"; M = 30; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 163; D = 44351; K = "Start Ruby"; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 164; D = 44494; K = "Start Ruby"; L = "Primera"; M = 200; N = 11000; O = 12000; P = 11500; S = 821 },
    @{ Row = 165; D = 44237; K = "Start Ruby"; L = "Primera"; M = 30; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 166; D = 44609; K = "Start Ruby"; L = "Primera"; M = 60; N = 14000; O = 14000; P = 14000; S = 1000 },
    @{ Row = 167; D = 44225; K = "Start Ruby"; L = "Primera"; M = 180; N = 14000; O = 15000; P = 14500; S = 1036 },
    @{ Row = 168; D = 44636; K = "Start Ruby"; L = "Primera"; M = 80; N = 12000; O = 13000; P = 12500; S = 893 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("D$i").Value = $r.D
    $ws.Range("K$i").Value = $r.K
    $ws.Range("L$i").Value = $r.L
    $ws.Range("M$i").Value = $r.M
    $ws.Range("N$i").Value = $r.N
    $ws.Range("O$i").Value = $r.O
    $ws.Range("P$i").Value = $r.P
    $ws.Range("S$i").Value = $r.S
}

Write-Output "Updated $($rows.Count) rows (132-168)."
